# "=sum(e5:e14) and corrected values in column e(rate)"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected Rate values (column E) for rows 6-14 (row 5 stays 2500).
# Columns F (Gross Salary), H (tax) and I (Net Pay) are formulas and
# recalculate automatically from the new rates.
$ws.Range("E6").Value = 5000
$ws.Range("E7").Value = 5000
$ws.Range("E8").Value = 5000
$ws.Range("E9").Value = 6500
$ws.Range("E10").Value = 4500
$ws.Range("E11").Value = 1200
$ws.Range("E12").Value = 3200
$ws.Range("E13").Value = 4000
$ws.Range("E14").Value = 2600

# New TOTAL row formula summing the Rate column.
$ws.Range("E15").Formula = "=SUM(E5:E14)"

# Leave the selection on the newly added total cell.
$ws.Range("E15").Select()
